function main() {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    $ws.Range("A2").Value = "7VlJ2"
    $ws.Range("A3").Value = "2DtB3"
    $ws.Range("A4").Value = "3EtC2"
    $ws.Range("A5").Value = "5SbY2"
    $ws.Range("A6").Value = "9YmO4"
}

main
